# Implemented functionality for dynamic template generation:
# add a new "template_to_send" column (D) that maps each row to the
# email template file that should be used when sending to that recruiter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("D1").Value = "template_to_send"

# Per-row template assignment
$ws.Range("D2").Value = "test_template.docx"
$ws.Range("D3").Value = "test_template_2.txt"
$ws.Range("D4").Value = "test_template.docx"

# Update the used range / selection to reflect the new column
$ws.Columns.Item(4).ColumnWidth = 15.2
$ws.Range("D4").Select()
